$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.379281
$ws.Range("H2").Value = 22.137843
$ws.Range("I2").Value = 0.2744121884499962
$ws.Range("J2").Value = 0.2744121884499961
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.701354
$ws.Range("N2").Value = 8.104061999999999
$ws.Range("O2").Value = 0.02221077311549548
$ws.Range("P2").Value = 0.02221077311549548
$ws.Range("Q2").Value = 19.934050246474
$ws.Range("R2").Value = 179.406452218266
$ws.Range("S2").Value = 0.006094906857789454
$ws.Range("T2").Value = 0.006094906857789453
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.379281
$ws.Range("H3").Value = 22.137843
$ws.Range("I3").Value = 0.2744121884499962
$ws.Range("J3").Value = 0.2744121884499961
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 88.14978533333333
$ws.Range("N3").Value = 264.449356
$ws.Range("O3").Value = 0.7247753838328104
$ws.Range("P3").Value = 0.7247753838328105
$ws.Range("Q3").Value = 650.4820360643453
$ws.Range("R3").Value = 5854.338324579107
$ws.Range("S3").Value = 0.1988871992122475
$ws.Range("T3").Value = 0.1988871992122475
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 7.379281
$ws.Range("H4").Value = 22.137843
$ws.Range("I4").Value = 0.2744121884499962
$ws.Range("J4").Value = 0.2744121884499961
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.24063
$ws.Range("N4").Value = 0.72189
$ws.Range("O4").Value = 0.001978481285600361
$ws.Range("P4").Value = 0.001978481285600361
$ws.Range("Q4").Value = 1.77567638703
$ws.Range("R4").Value = 15.98108748327
$ws.Range("S4").Value = 0.000542919379388957
$ws.Range("T4").Value = 0.0005429193793889569
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 7.379281
$ws.Range("H5").Value = 22.137843
$ws.Range("I5").Value = 0.2744121884499962
$ws.Range("J5").Value = 0.2744121884499961
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.53182233333333
$ws.Range("N5").Value = 91.595467
$ws.Range("O5").Value = 0.2510353617660938
$ws.Range("P5").Value = 0.2510353617660938
$ws.Range("Q5").Value = 225.3028964397423
$ws.Range("R5").Value = 2027.726067957681
$ws.Range("S5").Value = 0.0688871630005703
$ws.Range("T5").Value = 0.06888716300057028
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 13.29805733333333
$ws.Range("H6").Value = 39.894172
$ws.Range("I6").Value = 0.4945128143207339
$ws.Range("J6").Value = 0.4945128143207338
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.701354
$ws.Range("N6").Value = 8.104061999999999
$ws.Range("O6").Value = 0.02221077311549548
$ws.Range("P6").Value = 0.02221077311549548
$ws.Range("Q6").Value = 35.92276036962933
$ws.Range("R6").Value = 323.304843326664
$ws.Range("S6").Value = 0.01098351192158297
$ws.Range("T6").Value = 0.01098351192158296
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 13.29805733333333
$ws.Range("H7").Value = 39.894172
$ws.Range("I7").Value = 0.4945128143207339
$ws.Range("J7").Value = 0.4945128143207338
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 88.14978533333333
$ws.Range("N7").Value = 264.449356
$ws.Range("O7").Value = 0.7247753838328104
$ws.Range("P7").Value = 0.7247753838328105
$ws.Range("Q7").Value = 1172.220899283692
$ws.Range("R7").Value = 10549.98809355323
$ws.Range("S7").Value = 0.3584107148095532
$ws.Range("T7").Value = 0.3584107148095532
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 13.29805733333333
$ws.Range("H8").Value = 39.894172
$ws.Range("I8").Value = 0.4945128143207339
$ws.Range("J8").Value = 0.4945128143207338
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.24063
$ws.Range("N8").Value = 0.72189
$ws.Range("O8").Value = 0.001978481285600361
$ws.Range("P8").Value = 0.001978481285600361
$ws.Range("Q8").Value = 3.19991153612
$ws.Range("R8").Value = 28.79920382508
$ws.Range("S8").Value = 0.0009783843486231383
$ws.Range("T8").Value = 0.0009783843486231381
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 13.29805733333333
$ws.Range("H9").Value = 39.894172
$ws.Range("I9").Value = 0.4945128143207339
$ws.Range("J9").Value = 0.4945128143207338
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.53182233333333
$ws.Range("N9").Value = 91.595467
$ws.Range("O9").Value = 0.2510353617660938
$ws.Range("P9").Value = 0.2510353617660938
$ws.Range("Q9").Value = 406.0139238798138
$ws.Range("R9").Value = 3654.125314918324
$ws.Range("S9").Value = 0.1241402032409746
$ws.Range("T9").Value = 0.1241402032409746
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.673314666666667
$ws.Range("H10").Value = 5.019944000000001
$ws.Range("I10").Value = 0.06222529534320158
$ws.Range("J10").Value = 0.06222529534320156
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.701354
$ws.Range("N10").Value = 8.104061999999999
$ws.Range("O10").Value = 0.02221077311549548
$ws.Range("P10").Value = 0.02221077311549548
$ws.Range("Q10").Value = 4.520215268058667
$ws.Range("R10").Value = 40.681937412528
$ws.Range("S10").Value = 0.001382071916912548
$ws.Range("T10").Value = 0.001382071916912548
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.673314666666667
$ws.Range("H11").Value = 5.019944000000001
$ws.Range("I11").Value = 0.06222529534320158
$ws.Range("J11").Value = 0.06222529534320156
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 88.14978533333333
$ws.Range("N11").Value = 264.449356
$ws.Range("O11").Value = 0.7247753838328104
$ws.Range("P11").Value = 0.7247753838328105
$ws.Range("Q11").Value = 147.5023286617849
$ws.Range("R11").Value = 1327.520957956064
$ws.Range("S11").Value = 0.04509936231647891
$ws.Range("T11").Value = 0.04509936231647891
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.673314666666667
$ws.Range("H12").Value = 5.019944000000001
$ws.Range("I12").Value = 0.06222529534320158
$ws.Range("J12").Value = 0.06222529534320156
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.24063
$ws.Range("N12").Value = 0.72189
$ws.Range("O12").Value = 0.001978481285600361
$ws.Range("P12").Value = 0.001978481285600361
$ws.Range("Q12").Value = 0.4026497082400001
$ws.Range("R12").Value = 3.623847374160001
$ws.Range("S12").Value = 0.0001231115823274796
$ws.Range("T12").Value = 0.0001231115823274796
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.673314666666667
$ws.Range("H13").Value = 5.019944000000001
$ws.Range("I13").Value = 0.06222529534320158
$ws.Range("J13").Value = 0.06222529534320156
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 30.53182233333333
$ws.Range("N13").Value = 91.595467
$ws.Range("O13").Value = 0.2510353617660938
$ws.Range("P13").Value = 0.2510353617660938
$ws.Range("Q13").Value = 51.08934611042756
$ws.Range("R13").Value = 459.804114993848
$ws.Range("S13").Value = 0.01562074952748264
$ws.Range("T13").Value = 0.01562074952748264
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.540576
$ws.Range("H14").Value = 13.621728
$ws.Range("I14").Value = 0.1688497018860685
$ws.Range("J14").Value = 0.1688497018860685
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 2.701354
$ws.Range("N14").Value = 8.104061999999999
$ws.Range("O14").Value = 0.02221077311549548
$ws.Range("P14").Value = 0.02221077311549548
$ws.Range("Q14").Value = 12.265703139904
$ws.Range("R14").Value = 110.391328259136
$ws.Range("S14").Value = 0.003750282419210517
$ws.Range("T14").Value = 0.003750282419210517
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.540576
$ws.Range("H15").Value = 13.621728
$ws.Range("I15").Value = 0.1688497018860685
$ws.Range("J15").Value = 0.1688497018860685
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 88.14978533333333
$ws.Range("N15").Value = 264.449356
$ws.Range("O15").Value = 0.7247753838328104
$ws.Range("P15").Value = 0.7247753838328105
$ws.Range("Q15").Value = 400.2507996896853
$ws.Range("R15").Value = 3602.257197207167
$ws.Range("S15").Value = 0.1223781074945309
$ws.Range("T15").Value = 0.1223781074945309
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.540576
$ws.Range("H16").Value = 13.621728
$ws.Range("I16").Value = 0.1688497018860685
$ws.Range("J16").Value = 0.1688497018860685
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.24063
$ws.Range("N16").Value = 0.72189
$ws.Range("O16").Value = 0.001978481285600361
$ws.Range("P16").Value = 0.001978481285600361
$ws.Range("Q16").Value = 1.09259880288
$ws.Range("R16").Value = 9.83338922592
$ws.Range("S16").Value = 0.0003340659752607866
$ws.Range("T16").Value = 0.0003340659752607866
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 4.540576
$ws.Range("H17").Value = 13.621728
$ws.Range("I17").Value = 0.1688497018860685
$ws.Range("J17").Value = 0.1688497018860685
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 30.53182233333333
$ws.Range("N17").Value = 91.595467
$ws.Range("O17").Value = 0.2510353617660938
$ws.Range("P17").Value = 0.2510353617660938
$ws.Range("Q17").Value = 138.6320597229973
$ws.Range("R17").Value = 1247.688537506976
$ws.Range("S17").Value = 0.0423872459970663
$ws.Range("T17").Value = 0.04238724599706629
